$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'62.808.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.52%  "

# Row 3
$ws.Range("D3").Value = "'2.678.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.13%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "'553.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.01%  "

# Row 6
$ws.Range("D6").Value = "'157.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.93%  "

# Row 7
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("E8").Value = "  -1.05%  "

# Row 9
$ws.Range("D9").Value = "'0.105"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.50%  "

# Row 10
$ws.Range("D10").Value = "'0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.73%  "

# Row 11
$ws.Range("D11").Value = "'0.366"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.55%  "

# Row 12
$ws.Range("D12").Value = "'5.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.39%  "

# Row 13
$ws.Range("D13").Value = "'3.154.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.00%  "

# Row 14
$ws.Range("D14").Value = "'26.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.22%  "

# Row 15
$ws.Range("D15").Value = "'62.737.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.35%  "

# Row 16
$ws.Range("D16").Value = "'0.0000145"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.97%  "

# Row 17
$ws.Range("D17").Value = "'2.683.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.02%  "

# Row 18
$ws.Range("D18").Value = "'11.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.52%  "

# Row 19
$ws.Range("D19").Value = "'4.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.82%  "

# Row 20
$ws.Range("D20").Value = "'344.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.67%  "

# Row 21
$ws.Range("D21").Value = "'6.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.14%  "

# Row 22
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "

# Row 23
$ws.Range("D23").Value = "'0.510"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.15%  "

# Row 24
$ws.Range("D24").Value = "'63.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.63%  "

# Row 25
$ws.Range("D25").Value = "'0.169"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.05%  "

# Row 26
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.14%  "

# Row 27
$ws.Range("D27").Value = "'8.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.28%  "

# Row 28
$ws.Range("D28").Value = "'0.0₃0849"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.21%  "

# Row 29
$ws.Range("D29").Value = "'1.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.56%  "

# Row 30
$ws.Range("D30").Value = "'7.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.19%  "

# Row 31
$ws.Range("D31").Value = "'1.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.68%  "

# Row 32
$ws.Range("D32").Value = "'164.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.32%  "

# Row 33
$ws.Range("D33").Value = "'4.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.68%  "

# Row 35
$ws.Range("D35").Value = "'1.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.06%  "

# Row 36
$ws.Range("D36").Value = "'19.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.33%  "

# Row 37
$ws.Range("D37").Value = "'1.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.74%  "

# Row 38
$ws.Range("D38").Value = "'340.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.77%  "

# Row 39
$ws.Range("B39").Value = "SuiNetwork"
$ws.Range("C39").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D39").Value = "'0.934"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.39%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'6.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.80%  "

# Row 41
$ws.Range("D41").Value = "'3.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.41%  "

# Row 42
$ws.Range("D42").Value = "'38.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.51%  "

# Row 43
$ws.Range("D43").Value = "'20.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.91%  "

# Row 44
$ws.Range("D44").Value = "'0.617"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.25%  "

# Row 45
$ws.Range("D45").Value = "'20.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.66%  "

# Row 46
$ws.Range("D46").Value = "'0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.13%  "

# Row 47
$ws.Range("D47").Value = "'0.0553"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.13%  "

# Row 48
$ws.Range("E48").Value = "  -0.56%  "

# Row 49
$ws.Range("D49").Value = "'0.0966"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.40%  "

# Row 50
$ws.Range("D50").Value = "'128.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.11%  "

# Row 51
$ws.Range("D51").Value = "'0.0240"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.05%  "
